$wb = $excel.ActiveWorkbook

# Sheet ALC, row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3677.1428
$ws.Range("L40").Value = 4596.6665
$ws.Range("J40").Value = 4596.6665
$ws.Range("N40").Value = -4946.6665

# Sheet ALC, row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K62").Value = 4665.154
$ws.Range("H62").Value = 8912.333000000001
$ws.Range("I62").Value = 4665.154
$ws.Range("M62").Value = -4041.154

# Sheet ALC, row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 8912.333000000001
$ws.Range("K65").Value = 23325.77
$ws.Range("M65").Value = -20205.77
$ws.Range("I65").Value = 4665.154

# Sheet ALC, row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("L111").Value = 8000.000100000001
$ws.Range("H111").Value = 2969.6667
$ws.Range("K111").Value = 9363.500100000001
$ws.Range("M111").Value = -6296.500100000001
$ws.Range("J111").Value = 2666.6667
$ws.Range("I111").Value = 3121.1667
$ws.Range("N111").Value = -14134.0001

# Sheet ALC, row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I116").Value = 5464
$ws.Range("K116").Value = 5464
$ws.Range("H116").Value = 6138.6665
$ws.Range("M116").Value = -2022

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2331.6667
$ws.Range("J132").Value = 6969
$ws.Range("I132").Value = 2000.4286
$ws.Range("M132").Value = -3471.2858
$ws.Range("K132").Value = 6001.2858
$ws.Range("N132").Value = -25967
$ws.Range("L132").Value = 20907

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I138").Value = 1130
$ws.Range("J138").Value = 2929.0605
$ws.Range("N138").Value = -19067.1815
$ws.Range("M138").Value = 1750
$ws.Range("K138").Value = 3390
$ws.Range("L138").Value = 8787.181500000001
$ws.Range("H138").Value = 2294.0981

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4854.515
$ws.Range("J132").Value = 6777.3335
$ws.Range("I132").Value = 3252.1667
$ws.Range("M132").Value = -7226.500100000001
$ws.Range("K132").Value = 9756.500100000001
$ws.Range("N132").Value = -25392.0005
$ws.Range("L132").Value = 20332.0005

# Sheet BSM, row 75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M75").Value = -8063.5
$ws.Range("I75").Value = 8999.5
$ws.Range("L75").Value = 39080
$ws.Range("H75").Value = 30485.572
$ws.Range("J75").Value = 39080
$ws.Range("K75").Value = 8999.5
$ws.Range("N75").Value = -40952

# Sheet BSM, row 78
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I78").Value = 8999.5
$ws.Range("N78").Value = -126600
$ws.Range("J78").Value = 39080
$ws.Range("L78").Value = 117240
$ws.Range("H78").Value = 30485.572
$ws.Range("M78").Value = -22318.5
$ws.Range("K78").Value = 26998.5

# Sheet BSM, row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M86").Value = -1128.8
$ws.Range("I86").Value = 2251.8
$ws.Range("K86").Value = 2251.8
$ws.Range("H86").Value = 2251.8

# Sheet BSM, row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I89").Value = 2251.8
$ws.Range("K89").Value = 11259
$ws.Range("M89").Value = -5643
$ws.Range("H89").Value = 2251.8

# Sheet BSM, row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M94").Value = -977.8
$ws.Range("K94").Value = 1428.8
$ws.Range("H94").Value = 1445.174
$ws.Range("I94").Value = 1428.8

# Sheet BSM, row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2751.6
$ws.Range("J105").Value = 2628.5334
$ws.Range("N105").Value = -6122.5334
$ws.Range("M105").Value = -1373.8
$ws.Range("I105").Value = 3120.8
$ws.Range("K105").Value = 3120.8
$ws.Range("L105").Value = 2628.5334

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I134").Value = 1420.2858
$ws.Range("H134").Value = 436125.66
$ws.Range("M134").Value = -1725.857400000001
$ws.Range("K134").Value = 4260.857400000001

# Sheet CRP, row 5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N5").Value = -17650.857
$ws.Range("H5").Value = 15499.375
$ws.Range("J5").Value = 17426.857
$ws.Range("L5").Value = 17426.857

# Sheet CRP, row 53
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L53").Value = 51343.25
$ws.Range("J53").Value = 51343.25
$ws.Range("H53").Value = 51343.25
$ws.Range("N53").Value = -52557.25

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L58").Value = 3920.75
$ws.Range("N58").Value = -4326.75
$ws.Range("K58").Value = 1531
$ws.Range("H58").Value = 2034.1052
$ws.Range("I58").Value = 1531
$ws.Range("M58").Value = -1328
$ws.Range("J58").Value = 3920.75

# Sheet CRP, row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M94").Value = -6054
$ws.Range("J94").Value = 4534.636
$ws.Range("K94").Value = 6505
$ws.Range("H94").Value = 4837.769
$ws.Range("I94").Value = 6505
$ws.Range("N94").Value = -5436.636
$ws.Range("L94").Value = 4534.636

# Sheet CRP, row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1156.4706
$ws.Range("M105").Value = 550.9231
$ws.Range("I105").Value = 1196.0769
$ws.Range("K105").Value = 1196.0769

# Sheet CRP, row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M107").Value = 870.1666
$ws.Range("J107").Value = 2938
$ws.Range("I107").Value = 1049.8334
$ws.Range("N107").Value = -6778
$ws.Range("L107").Value = 2938
$ws.Range("H107").Value = 2271.5881
$ws.Range("K107").Value = 1049.8334

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2013.125
$ws.Range("J132").Value = 0
$ws.Range("I132").Value = 2013.125
$ws.Range("M132").Value = -3509.375
$ws.Range("K132").Value = 6039.375
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I134").Value = 1805.75
$ws.Range("H134").Value = 2943.6316
$ws.Range("M134").Value = -2882.25
$ws.Range("K134").Value = 5417.25

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I136").Value = 1531
$ws.Range("N136").Value = -16862.25
$ws.Range("J136").Value = 3920.75
$ws.Range("M136").Value = -2043
$ws.Range("H136").Value = 2034.1052
$ws.Range("L136").Value = 11762.25
$ws.Range("K136").Value = 4593

# Sheet CUL, row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2999.5
$ws.Range("K68").Value = 2997
$ws.Range("I68").Value = 999
$ws.Range("M68").Value = -2186

# Sheet CUL, row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2999.5
$ws.Range("K71").Value = 8991
$ws.Range("I71").Value = 999
$ws.Range("M71").Value = -4935

# Sheet GSM, row 95
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 125041420
$ws.Range("J95").Value = 125041420
$ws.Range("N95").Value = -125046912
$ws.Range("L95").Value = 125041420

# Sheet GSM, row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N97").Value = -1891
$ws.Range("H97").Value = 1734.3846
$ws.Range("K97").Value = 2450.4285
$ws.Range("I97").Value = 2450.4285
$ws.Range("M97").Value = -1954.4285
$ws.Range("J97").Value = 899
$ws.Range("L97").Value = 899

# Sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1572.775
$ws.Range("M102").Value = 432.3214
$ws.Range("I102").Value = 1189.6786
$ws.Range("N102").Value = -5710.6667
$ws.Range("L102").Value = 2466.6667
$ws.Range("J102").Value = 2466.6667
$ws.Range("K102").Value = 1189.6786

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 24396104
$ws.Range("J132").Value = 17924
$ws.Range("I132").Value = 29415140
$ws.Range("M132").Value = -88242890
$ws.Range("K132").Value = 88245420
$ws.Range("N132").Value = -58832
$ws.Range("L132").Value = 53772

# Sheet LTW, row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N46").Value = -11918.857
$ws.Range("J46").Value = 11542.857
$ws.Range("H46").Value = 4764.52
$ws.Range("L46").Value = 11542.857

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M122").Value = -11456.059
$ws.Range("H122").Value = 5276.2383
$ws.Range("K122").Value = 13906.059
$ws.Range("I122").Value = 4635.353

# Sheet WVR, row 39
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 27198.334
$ws.Range("I39").Value = 25545
$ws.Range("K39").Value = 25545
$ws.Range("M39").Value = -25132

# Sheet WVR, row 45
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N45").Value = -12982
$ws.Range("L45").Value = 12000
$ws.Range("H45").Value = 12000
$ws.Range("J45").Value = 12000

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2111.7334
$ws.Range("I132").Value = 1548.2858
$ws.Range("M132").Value = -2114.857400000001
$ws.Range("K132").Value = 4644.857400000001
